# Scheduled runner update: refresh market-price-derived columns (H:N)
# across the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1680.5186
$ws.Range("I40").Value = 1569.579
$ws.Range("J40").Value = 1944
$ws.Range("K40").Value = 1569.579
$ws.Range("L40").Value = 1944
$ws.Range("M40").Value = -1394.579
$ws.Range("N40").Value = -2294
$ws.Range("H86").Value = 1354.6
$ws.Range("I86").Value = 1354.6
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1354.6
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -231.5999999999999
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1354.6
$ws.Range("I89").Value = 1354.6
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6773
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1157
$ws.Range("N89").ClearContents()
$ws.Range("H92").Value = 22224462
$ws.Range("I92").Value = 22224462
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 22224462
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -22223214
$ws.Range("N92").ClearContents()
$ws.Range("H98").Value = 659298
$ws.Range("I98").Value = 799047.5600000001
$ws.Range("J98").Value = 7133.3335
$ws.Range("K98").Value = 799047.5600000001
$ws.Range("L98").Value = 7133.3335
$ws.Range("M98").Value = -797549.5600000001
$ws.Range("N98").Value = -10129.3335
$ws.Range("H106").Value = 55557556
$ws.Range("I106").Value = 55557556
$ws.Range("K106").Value = 55557556
$ws.Range("M106").Value = -55556925
$ws.Range("H111").Value = 988.2308
$ws.Range("I111").Value = 634.7
$ws.Range("J111").Value = 2166.6667
$ws.Range("K111").Value = 1904.1
$ws.Range("L111").Value = 6500.000100000001
$ws.Range("M111").Value = 1162.9
$ws.Range("N111").Value = -12634.0001
$ws.Range("H112").Value = 8065611
$ws.Range("J112").Value = 8334458
$ws.Range("L112").Value = 25003374
$ws.Range("N112").Value = -25005590
$ws.Range("H122").Value = 659298
$ws.Range("I122").Value = 799047.5600000001
$ws.Range("J122").Value = 7133.3335
$ws.Range("K122").Value = 2397142.68
$ws.Range("L122").Value = 21400.0005
$ws.Range("M122").Value = -2394692.68
$ws.Range("N122").Value = -26300.0005
$ws.Range("H132").Value = 372906.22
$ws.Range("I132").Value = 675839.25
$ws.Range("J132").Value = 32106.562
$ws.Range("K132").Value = 2027517.75
$ws.Range("L132").Value = 96319.686
$ws.Range("M132").Value = -2024987.75
$ws.Range("N132").Value = -101379.686
$ws.Range("H138").Value = 1726.8
$ws.Range("I138").Value = 761.4091
$ws.Range("J138").Value = 1999.0897
$ws.Range("K138").Value = 2284.2273
$ws.Range("L138").Value = 5997.2691
$ws.Range("M138").Value = 2855.7727
$ws.Range("N138").Value = -16277.2691

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5650
$ws.Range("I2").Value = 7067.3687
$ws.Range("J2").Value = 1161.6666
$ws.Range("K2").Value = 7067.3687
$ws.Range("L2").Value = 1161.6666
$ws.Range("M2").Value = -6954.3687
$ws.Range("N2").Value = -1387.6666
$ws.Range("H45").Value = 2131.6155
$ws.Range("I45").Value = 2059.25
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 2059.25
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1682.25
$ws.Range("N45").Value = -3754
$ws.Range("H74").Value = 9463.467000000001
$ws.Range("I74").Value = 1684.909
$ws.Range("J74").Value = 30854.5
$ws.Range("K74").Value = 1684.909
$ws.Range("L74").Value = 30854.5
$ws.Range("M74").Value = -810.9090000000001
$ws.Range("N74").Value = -32602.5
$ws.Range("H77").Value = 9463.467000000001
$ws.Range("I77").Value = 1684.909
$ws.Range("J77").Value = 30854.5
$ws.Range("K77").Value = 8424.545
$ws.Range("L77").Value = 154272.5
$ws.Range("M77").Value = -4056.545
$ws.Range("N77").Value = -163008.5
$ws.Range("H97").Value = 47633920
$ws.Range("I97").Value = 66687212
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 66687212
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -66686716
$ws.Range("N97").Value = -1692
$ws.Range("H116").Value = 5650
$ws.Range("I116").Value = 7067.3687
$ws.Range("J116").Value = 1161.6666
$ws.Range("K116").Value = 7067.3687
$ws.Range("L116").Value = 1161.6666
$ws.Range("M116").Value = -4773.3687
$ws.Range("N116").Value = -5749.6666
$ws.Range("H132").Value = 2701.2903
$ws.Range("I132").Value = 2226.75
$ws.Range("J132").Value = 4328.2856
$ws.Range("K132").Value = 6680.25
$ws.Range("L132").Value = 12984.8568
$ws.Range("M132").Value = -4150.25
$ws.Range("N132").Value = -18044.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5650
$ws.Range("I3").Value = 7067.3687
$ws.Range("J3").Value = 1161.6666
$ws.Range("K3").Value = 7067.3687
$ws.Range("L3").Value = 1161.6666
$ws.Range("M3").Value = -6953.3687
$ws.Range("N3").Value = -1389.6666
$ws.Range("H20").Value = 1466.3572
$ws.Range("I20").Value = 1313.6875
$ws.Range("J20").Value = 1669.9166
$ws.Range("K20").Value = 1313.6875
$ws.Range("L20").Value = 1669.9166
$ws.Range("M20").Value = -1066.6875
$ws.Range("N20").Value = -2163.9166
$ws.Range("H134").Value = 2693.9583
$ws.Range("I134").Value = 1935.05
$ws.Range("K134").Value = 5805.15
$ws.Range("M134").Value = -3270.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3288.2778
$ws.Range("I31").Value = 1379.6562
$ws.Range("K31").Value = 1379.6562
$ws.Range("M31").Value = -1084.6562
$ws.Range("H34").Value = 3288.2778
$ws.Range("I34").Value = 1379.6562
$ws.Range("K34").Value = 1379.6562
$ws.Range("M34").Value = -1177.6562
$ws.Range("H132").Value = 3446.611
$ws.Range("I132").Value = 1764.8889
$ws.Range("K132").Value = 5294.6667
$ws.Range("M132").Value = -2764.6667
$ws.Range("H134").Value = 2847.35
$ws.Range("I134").Value = 2000.9667
$ws.Range("J134").Value = 5386.5
$ws.Range("K134").Value = 6002.9001
$ws.Range("L134").Value = 16159.5
$ws.Range("M134").Value = -3467.9001
$ws.Range("N134").Value = -21229.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1189.7142
$ws.Range("J122").Value = 1932.6666
$ws.Range("L122").Value = 17393.9994
$ws.Range("N122").Value = -22293.9994
$ws.Range("H131").Value = 2574.8354
$ws.Range("J131").Value = 2664.8948
$ws.Range("L131").Value = 7994.6844
$ws.Range("N131").Value = -18074.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 13000
$ws.Range("J96").Value = 13000
$ws.Range("L96").Value = 13000
$ws.Range("N96").Value = -18492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1825.8572
$ws.Range("I46").Value = 1293.6666
$ws.Range("J46").Value = 2225
$ws.Range("K46").Value = 1293.6666
$ws.Range("L46").Value = 2225
$ws.Range("M46").Value = -1105.6666
$ws.Range("N46").Value = -2601
$ws.Range("H132").Value = 4797.778
$ws.Range("I132").Value = 2995.25
$ws.Range("J132").Value = 6239.8
$ws.Range("K132").Value = 8985.75
$ws.Range("L132").Value = 18719.4
$ws.Range("M132").Value = -6455.75
$ws.Range("N132").Value = -23779.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 59999
$ws.Range("J95").Value = 59999
$ws.Range("L95").Value = 59999
$ws.Range("N95").Value = -65491
$ws.Range("H107").Value = 622.80646
$ws.Range("I107").Value = 628.7143
$ws.Range("J107").Value = 610.4
$ws.Range("K107").Value = 1886.1429
$ws.Range("L107").Value = 1831.2
$ws.Range("M107").Value = 33.85710000000017
$ws.Range("N107").Value = -5671.2
$ws.Range("H132").Value = 21744432
$ws.Range("I132").Value = 31255544
$ws.Range("J132").Value = 4749.2856
$ws.Range("K132").Value = 93766632
$ws.Range("L132").Value = 14247.8568
$ws.Range("M132").Value = -93764102
$ws.Range("N132").Value = -19307.8568
$ws.Range("H136").Value = 37040290
$ws.Range("I136").Value = 111112380
$ws.Range("J136").Value = 4242.1665
$ws.Range("K136").Value = 333337140
$ws.Range("L136").Value = 12726.4995
$ws.Range("M136").Value = -333334590
$ws.Range("N136").Value = -17826.4995
